# Add a "Save" column (H) to the s_vals sheet.
#
# Header cell H1 gets the same look as the other header cells (e.g. G1):
# bold, centered, thin-bordered. We copy G1's formatting via copy/paste-special
# so the new cell reuses the existing style entry instead of Excel minting a
# brand new (slightly different) cellXf/font.
#
# Data cells H2 / H3 are plain numeric zeros, matching the unstyled numeric
# cells elsewhere in the sheet (e.g. F2/F3, G2/G3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "sum" header (G1) onto the new H1 cell.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)   # xlPasteFormats

# Now set the header text and the two data values.
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
